# Fiducial Worksheet Transforms - apply updated surface measurement values
# and refresh the saved cursor/selection position, per the latest
# procedures/measurements/reports update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "180 to 0" (column C) fiducial measurements.
$ws.Range("C2").Value = -2.27
$ws.Range("C3").Value = 0.29
$ws.Range("C4").Value = -179.65
$ws.Range("C5").Value = 1

# The workbook was last left with the cursor on D6.
$ws.Range("D6").Select() | Out-Null
